# Weekly price-sheet update: a new week's record is inserted as row 20,
# pushing the previously existing rows 20..53 down to 21..54.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 20 (shifts rows 20-53 -> 21-54,
# and extends the used range from R53 to R54).
$ws.Rows(20).Insert()

# Populate the new row 20 with the new weekly record.
$ws.Cells.Item(20, 1).Value = 10
$ws.Cells.Item(20, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(20, 3).Value = "La Araucanía"
$ws.Cells.Item(20, 4).Value = 44622
$ws.Cells.Item(20, 5).Value = 9
$ws.Cells.Item(20, 6).Value = 100114002
$ws.Cells.Item(20, 7).Value = "Camote"
$ws.Cells.Item(20, 8).Value = "Sin especificar"
$ws.Cells.Item(20, 9).Value = "Primera"
$ws.Cells.Item(20, 10).Value = 50
$ws.Cells.Item(20, 11).Value = 18000
$ws.Cells.Item(20, 12).Value = 18000
$ws.Cells.Item(20, 13).Value = 18000
$ws.Cells.Item(20, 14).Value = "$/malla 20 kilos"
$ws.Cells.Item(20, 15).Value = "Perú"
$ws.Cells.Item(20, 16).Value = 900
$ws.Cells.Item(20, 17).Value = 20
$ws.Cells.Item(20, 18).Value = "Hortaliza"
